# Auto-generated: update market price snapshot columns (H-N) per scheduled runner refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 468.875
$ws.Range("J19").Value = 357.85715
$ws.Range("L19").Value = 357.85715
$ws.Range("N19").Value = -707.85715
# Row 33
$ws.Range("H33").Value = 1944
$ws.Range("I33").Value = 1729.8
$ws.Range("J33").Value = 2479.5
$ws.Range("K33").Value = 1729.8
$ws.Range("L33").Value = 2479.5
$ws.Range("M33").Value = -1500.8
$ws.Range("N33").Value = -2937.5
# Row 41
$ws.Range("H41").Value = 1642.8182
$ws.Range("I41").Value = 1125.6666
$ws.Range("J41").Value = 2263.4
$ws.Range("K41").Value = 1125.6666
$ws.Range("L41").Value = 2263.4
$ws.Range("M41").Value = -685.6666
$ws.Range("N41").Value = -3143.4
# Row 62
$ws.Range("H62").Value = 4421.4
$ws.Range("I62").Value = 4421.4
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 4421.4
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -3797.4
$ws.Range("N62").ClearContents()
# Row 65
$ws.Range("H65").Value = 4421.4
$ws.Range("I65").Value = 4421.4
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 22107
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -18987
$ws.Range("N65").ClearContents()
# Row 112
$ws.Range("H112").Value = 1878.0605
$ws.Range("J112").Value = 1878.0605
$ws.Range("L112").Value = 5634.181500000001
$ws.Range("N112").Value = -7850.181500000001
# Row 138
$ws.Range("H138").Value = 1817.1794
$ws.Range("I138").Value = 1110.3226
$ws.Range("K138").Value = 3330.9678
$ws.Range("M138").Value = 1809.0322

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 9254.829
$ws.Range("I32").Value = 7543.8154
$ws.Range("J32").Value = 19365.363
$ws.Range("K32").Value = 7543.8154
$ws.Range("L32").Value = 19365.363
$ws.Range("M32").Value = -7256.8154
$ws.Range("N32").Value = -19939.363
# Row 88
$ws.Range("H88").Value = 2280.8125
$ws.Range("J88").Value = 2063.5
$ws.Range("L88").Value = 2063.5
$ws.Range("N88").Value = -2875.5
# Row 91
$ws.Range("H91").Value = 2280.8125
$ws.Range("J91").Value = 2063.5
$ws.Range("L91").Value = 2063.5
$ws.Range("N91").Value = -4871.5
# Row 97
$ws.Range("H97").Value = 1485.5834
$ws.Range("I97").Value = 1166.091
$ws.Range("K97").Value = 1166.091
$ws.Range("M97").Value = -670.0909999999999
# Row 110
$ws.Range("H110").Value = 1923.7391
$ws.Range("I110").Value = 2012.1904
$ws.Range("K110").Value = 2012.1904
$ws.Range("M110").Value = 32.80960000000005

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 29727
$ws.Range("I99").Value = 31946
$ws.Range("J99").Value = 880
$ws.Range("K99").Value = 31946
$ws.Range("L99").Value = 880
$ws.Range("M99").Value = -30448
$ws.Range("N99").Value = -3876
# Row 100
$ws.Range("H100").Value = 24065
$ws.Range("J100").Value = 24065
$ws.Range("L100").Value = 24065
$ws.Range("N100").Value = -26229
# Row 105
$ws.Range("H105").Value = 1802.4166
$ws.Range("I105").Value = 1875.3636
$ws.Range("K105").Value = 1875.3636
$ws.Range("M105").Value = -128.3635999999999

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 17589574
$ws.Range("I16").Value = 650.5
$ws.Range("K16").Value = 650.5
$ws.Range("M16").Value = -363.5
# Row 18
$ws.Range("H18").Value = 99000
$ws.Range("J18").Value = 99000
$ws.Range("L18").Value = 99000
$ws.Range("N18").Value = -99460
# Row 32
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
# Row 43
$ws.Range("H43").Value = 17240.875
$ws.Range("J43").Value = 17240.875
$ws.Range("L43").Value = 17240.875
$ws.Range("N43").Value = -17608.875
# Row 55
$ws.Range("H55").Value = 7810.4
$ws.Range("I55").Value = 7810.4
$ws.Range("K55").Value = 7810.4
$ws.Range("M55").Value = -7495.4
# Row 58
$ws.Range("H58").Value = 4082.3333
$ws.Range("J58").Value = 6665
$ws.Range("L58").Value = 6665
$ws.Range("N58").Value = -7071
# Row 101
$ws.Range("H101").Value = 17240.875
$ws.Range("J101").Value = 17240.875
$ws.Range("L101").Value = 17240.875
$ws.Range("N101").Value = -23730.875
# Row 106
$ws.Range("H106").Value = 37333.332
$ws.Range("J106").Value = 37333.332
$ws.Range("L106").Value = 37333.332
$ws.Range("N106").Value = -39857.332
# Row 113
$ws.Range("H113").Value = 17589574
$ws.Range("I113").Value = 650.5
$ws.Range("K113").Value = 650.5
$ws.Range("M113").Value = 1519.5
# Row 136
$ws.Range("H136").Value = 4082.3333
$ws.Range("J136").Value = 6665
$ws.Range("L136").Value = 19995
$ws.Range("N136").Value = -25095

$ws = $wb.Worksheets.Item("CUL")
# Row 122
$ws.Range("H122").Value = 12347059
$ws.Range("I122").Value = 1092.2858
$ws.Range("J122").Value = 16668147
$ws.Range("K122").Value = 9830.572200000001
$ws.Range("L122").Value = 150013323
$ws.Range("M122").Value = -7380.572200000001
$ws.Range("N122").Value = -150018223
# Row 133
$ws.Range("H133").Value = 1077.5

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 3333606.2
$ws.Range("I2").Value = 279.94446
$ws.Range("J2").Value = 8333595.5
$ws.Range("K2").Value = 279.94446
$ws.Range("L2").Value = 8333595.5
$ws.Range("M2").Value = -166.94446
$ws.Range("N2").Value = -8333821.5
# Row 70
$ws.Range("H70").Value = 165183.58
$ws.Range("J70").Value = 8868.75
$ws.Range("L70").Value = 8868.75
$ws.Range("N70").Value = -9408.75
# Row 73
$ws.Range("H73").Value = 165183.58
$ws.Range("J73").Value = 8868.75
$ws.Range("L73").Value = 8868.75
$ws.Range("N73").Value = -10740.75
# Row 80
$ws.Range("H80").Value = 95600.336
$ws.Range("I80").Value = 161044.14
$ws.Range("J80").Value = 3979
$ws.Range("K80").Value = 161044.14
$ws.Range("L80").Value = 3979
$ws.Range("M80").Value = -160046.14
$ws.Range("N80").Value = -5975
# Row 83
$ws.Range("H83").Value = 95600.336
$ws.Range("I83").Value = 161044.14
$ws.Range("J83").Value = 3979
$ws.Range("K83").Value = 805220.7000000001
$ws.Range("L83").Value = 19895
$ws.Range("M83").Value = -800228.7000000001
$ws.Range("N83").Value = -29879
# Row 97
$ws.Range("H97").Value = 374.5
$ws.Range("I97").Value = 374.5
$ws.Range("K97").Value = 374.5
$ws.Range("M97").Value = 121.5
# Row 105
$ws.Range("H105").Value = 69158.164
$ws.Range("J105").Value = 69158.164
$ws.Range("L105").Value = 69158.164
$ws.Range("N105").Value = -76146.164

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1092.6
$ws.Range("I22").Value = 662.375
$ws.Range("K22").Value = 662.375
$ws.Range("M22").Value = -367.375
# Row 27
$ws.Range("H27").Value = 1092.6
$ws.Range("I27").Value = 662.375
$ws.Range("K27").Value = 662.375
$ws.Range("M27").Value = -555.375
# Row 46
$ws.Range("H46").Value = 4316.8477
$ws.Range("I46").Value = 703.26086
$ws.Range("J46").Value = 7930.4346
$ws.Range("K46").Value = 703.26086
$ws.Range("L46").Value = 7930.4346
$ws.Range("M46").Value = -515.26086
$ws.Range("N46").Value = -8306.434600000001
# Row 47
$ws.Range("H47").Value = 34495
$ws.Range("J47").Value = 34495
$ws.Range("L47").Value = 34495
$ws.Range("N47").Value = -35475
# Row 52
$ws.Range("H52").Value = 34495
$ws.Range("J52").Value = 34495
$ws.Range("L52").Value = 34495
$ws.Range("N52").Value = -34961
# Row 93
$ws.Range("H93").Value = 592142.8
$ws.Range("I93").Value = 3887.7856
$ws.Range("J93").Value = 3337333
$ws.Range("K93").Value = 3887.7856
$ws.Range("L93").Value = 3337333
$ws.Range("M93").Value = -2639.7856
$ws.Range("N93").Value = -3339829

$ws = $wb.Worksheets.Item("WVR")
# Row 64
$ws.Range("H64").Value = 59886
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 59886
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 59886
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -60382
# Row 67
$ws.Range("H67").Value = 59886
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 59886
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 59886
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -61602
# Row 104
$ws.Range("H104").Value = 42185
$ws.Range("J104").Value = 42185
$ws.Range("L104").Value = 42185
$ws.Range("N104").Value = -49173
# Row 113
$ws.Range("H113").Value = 599.5714
$ws.Range("I113").Value = 566.1667
$ws.Range("K113").Value = 1698.5001
$ws.Range("M113").Value = 471.4999

